$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1831
$ws.Range("I18").Value = 1831
$ws.Range("K18").Value = 1831
$ws.Range("M18").Value = -1547
$ws.Range("H32").Value = 8000.5
$ws.Range("J32").Value = 8000
$ws.Range("L32").Value = 8000
$ws.Range("N32").Value = -8652
$ws.Range("H86").Value = 3339.1667
$ws.Range("H89").Value = 3339.1667
$ws.Range("H98").Value = 1586.1538
$ws.Range("I98").Value = 1466.1818
$ws.Range("K98").Value = 1466.1818
$ws.Range("M98").Value = 31.81819999999993
$ws.Range("H100").Value = 6696.849
$ws.Range("I100").Value = 1940.6111
$ws.Range("K100").Value = 1940.6111
$ws.Range("M100").Value = -1399.6111
$ws.Range("H116").Value = 4111.1113
$ws.Range("I116").Value = 2450
$ws.Range("J116").Value = 4941.6665
$ws.Range("K116").Value = 2450
$ws.Range("L116").Value = 4941.6665
$ws.Range("M116").Value = 992
$ws.Range("N116").Value = -11825.6665
$ws.Range("H122").Value = 1586.1538
$ws.Range("I122").Value = 1466.1818
$ws.Range("K122").Value = 4398.5454
$ws.Range("M122").Value = -1948.5454
$ws.Range("H137").Value = 5685900.5
$ws.Range("I137").Value = 7355000.5
$ws.Range("J137").Value = 10960
$ws.Range("K137").Value = 22065001.5
$ws.Range("L137").Value = 32880
$ws.Range("M137").Value = -22062451.5
$ws.Range("N137").Value = -37980
$ws.Range("H138").Value = 2598.7073
$ws.Range("J138").Value = 3978.25
$ws.Range("L138").Value = 11934.75
$ws.Range("N138").Value = -22214.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 29411.643
$ws.Range("I45").Value = 33316.168
$ws.Range("K45").Value = 33316.168
$ws.Range("M45").Value = -32939.168
$ws.Range("H74").Value = 299151.5
$ws.Range("I74").Value = 326811.3
$ws.Range("J74").Value = 13333.333
$ws.Range("K74").Value = 326811.3
$ws.Range("L74").Value = 13333.333
$ws.Range("M74").Value = -325937.3
$ws.Range("N74").Value = -15081.333
$ws.Range("H77").Value = 299151.5
$ws.Range("I77").Value = 326811.3
$ws.Range("J77").Value = 13333.333
$ws.Range("K77").Value = 1634056.5
$ws.Range("L77").Value = 66666.66500000001
$ws.Range("M77").Value = -1629688.5
$ws.Range("N77").Value = -75402.66500000001
$ws.Range("H102").Value = 1385.3
$ws.Range("I102").Value = 1385.3
$ws.Range("K102").Value = 1385.3
$ws.Range("M102").Value = 236.7
$ws.Range("H132").Value = 6286.696
$ws.Range("I132").Value = 5110.647
$ws.Range("J132").Value = 9618.833000000001
$ws.Range("K132").Value = 15331.941
$ws.Range("L132").Value = 28856.499
$ws.Range("M132").Value = -12801.941
$ws.Range("N132").Value = -33916.499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3517.8572
$ws.Range("J99").Value = 1299.5
$ws.Range("L99").Value = 1299.5
$ws.Range("N99").Value = -4295.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20836180
$ws.Range("I31").Value = 30304874
$ws.Range("J31").Value = 5056.2666
$ws.Range("K31").Value = 30304874
$ws.Range("L31").Value = 5056.2666
$ws.Range("M31").Value = -30304579
$ws.Range("N31").Value = -5646.2666
$ws.Range("H34").Value = 20836180
$ws.Range("I34").Value = 30304874
$ws.Range("J34").Value = 5056.2666
$ws.Range("K34").Value = 30304874
$ws.Range("L34").Value = 5056.2666
$ws.Range("M34").Value = -30304672
$ws.Range("N34").Value = -5460.2666
$ws.Range("H58").Value = 4594
$ws.Range("I58").Value = 3318.2
$ws.Range("K58").Value = 3318.2
$ws.Range("M58").Value = -3115.2
$ws.Range("H86").Value = 7343.0586
$ws.Range("I86").Value = 7359.5
$ws.Range("J86").Value = 7266.3335
$ws.Range("K86").Value = 7359.5
$ws.Range("L86").Value = 7266.3335
$ws.Range("M86").Value = -6236.5
$ws.Range("N86").Value = -9512.333500000001
$ws.Range("H89").Value = 7343.0586
$ws.Range("I89").Value = 7359.5
$ws.Range("J89").Value = 7266.3335
$ws.Range("K89").Value = 36797.5
$ws.Range("L89").Value = 36331.6675
$ws.Range("M89").Value = -31181.5
$ws.Range("N89").Value = -47563.6675
$ws.Range("H99").Value = 4546.2
$ws.Range("I99").Value = 4681.75
$ws.Range("J99").Value = 4004
$ws.Range("K99").Value = 4681.75
$ws.Range("L99").Value = 4004
$ws.Range("M99").Value = -3183.75
$ws.Range("N99").Value = -7000
$ws.Range("H126").Value = 4546.2
$ws.Range("I126").Value = 4681.75
$ws.Range("J126").Value = 4004
$ws.Range("K126").Value = 14045.25
$ws.Range("L126").Value = 12012
$ws.Range("M126").Value = -11575.25
$ws.Range("N126").Value = -16952
$ws.Range("H132").Value = 83105.59
$ws.Range("I132").Value = 86277.28999999999
$ws.Range("K132").Value = 258831.87
$ws.Range("M132").Value = -256301.87
$ws.Range("H136").Value = 4594
$ws.Range("I136").Value = 3318.2
$ws.Range("K136").Value = 9954.599999999999
$ws.Range("M136").Value = -7404.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 7960.125
$ws.Range("J107").Value = 7960.125
$ws.Range("L107").Value = 23880.375
$ws.Range("N107").Value = -27720.375
$ws.Range("H138").Value = 3473534.5
$ws.Range("I138").Value = 5000
$ws.Range("K138").Value = 15000
$ws.Range("M138").Value = -9860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6005.5835
$ws.Range("I80").Value = 4973.125
$ws.Range("K80").Value = 4973.125
$ws.Range("M80").Value = -3975.125
$ws.Range("H83").Value = 6005.5835
$ws.Range("I83").Value = 4973.125
$ws.Range("K83").Value = 24865.625
$ws.Range("M83").Value = -19873.625
$ws.Range("H95").Value = 33000
$ws.Range("J95").Value = 33000
$ws.Range("L95").Value = 33000
$ws.Range("N95").Value = -38492
$ws.Range("H97").Value = 950
$ws.Range("I97").Value = 900
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 900
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -404
$ws.Range("N97").Value = -1992
$ws.Range("H126").Value = 3574
$ws.Range("I126").Value = 2586.75
$ws.Range("J126").Value = 4451.5557
$ws.Range("K126").Value = 7760.25
$ws.Range("L126").Value = 13354.6671
$ws.Range("M126").Value = -5290.25
$ws.Range("N126").Value = -18294.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3477.8684
$ws.Range("I136").Value = 2348.6333
$ws.Range("K136").Value = 7045.8999
$ws.Range("M136").Value = -4495.8999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4788.685
$ws.Range("I132").Value = 3329.122
$ws.Range("J132").Value = 9391.923000000001
$ws.Range("K132").Value = 9987.366
$ws.Range("L132").Value = 28175.769
$ws.Range("M132").Value = -7457.366
$ws.Range("N132").Value = -33235.769

Write-Output "Applied all cell updates"